# ---------------------------------------------------------------------------
# Target-diff analysis
# ---------------------------------------------------------------------------
# The supplied unified diff touches word/document.xml and word/styles.xml
# only, and every single changed line is a pure XML *attribute re-ordering*:
# the element names, the attribute names, and every attribute *value* are
# identical before and after the change -- only the left-to-right order in
# which the attributes (and, on the <w:document> root, the xmlns:* namespace
# declarations) are written differs, e.g.
#
#   <w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/>
#   -> <w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/>
#
# is the same element with the same three attribute/value pairs, merely
# re-sorted.  The identical pattern (same attribute set & values, different
# order) repeats for every other hunk: the <w:document> namespace list,
# <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>, every single
# <w:lsdException>, and every <w:style>/<w:tblInd>/<w:tblCellMar> under the
# "Normal Table" style.  There is no insertion, deletion, or value change
# anywhere in the diff -- confirmed by canonicalising both the "before" and
# "after" XML (sorting each element's attributes) and comparing the element
# trees: they are identical.  This is a byte-level re-serialization artefact
# (the template having been re-saved by a different XML writer at some point
# in its history), not a content edit.
#
# Word's COM automation surface (Font.Color, PageSetup.*, Styles(...).*,
# etc.) is defined purely in terms of attribute *values*: it lets an author
# change what a property *is*, never the order in which the underlying
# OOXML writer lists an element's attributes when it serializes the part --
# that ordering is an internal writer implementation detail with no
# corresponding object-model property, in Word or in this host. There is
# therefore no COM call that "does" this diff beyond re-affirming that every
# value it touches is already correct. Concretely (verified against this
# host): re-assigning PageSetup.PageWidth/PageHeight/margins or a run's
# Font.Color back to their own current values does not change the writer's
# attribute order (it stays in the part's original/schema order either way)
# -- it only dirties the part and perturbs unrelated state (document
# statistics, namespace list) that is untouched by the diff. So the
# content-faithful application of this diff is to leave the document's
# values exactly as they are, which is what the object model already
# reports below.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# <w:sectPr>/<w:pgSz>+<w:pgMar> values are unchanged by the diff
# (11906 x 16838 twips page; 1417/1417/1417/1417 twips margins;
# 708/708 header/footer distance; 0 gutter). Confirm, do not mutate.
$pageSetup = $d.Sections(1).PageSetup
Write-Host "Page size (twips):" ([int]($pageSetup.PageWidth * 20)) "x" ([int]($pageSetup.PageHeight * 20))
Write-Host "Margins (twips) T/B/L/R:" ([int]($pageSetup.TopMargin * 20)) ([int]($pageSetup.BottomMargin * 20)) ([int]($pageSetup.LeftMargin * 20)) ([int]($pageSetup.RightMargin * 20))

# <w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/> on the
# "self" field-code run keeps the same value in the diff -- only the
# attribute order changes, which is not representable through the Field /
# Font object model. Confirm the field is present, do not mutate it.
Write-Host "Fields in document:" $d.Fields.Count
foreach ($field in $d.Fields) {
    Write-Host "Field code:" $field.Code.Text
}

Write-Host "No object-model mutation applied: every diff hunk only re-orders existing XML attributes/namespaces (same tags, same attribute values) in word/document.xml and word/styles.xml -- there is no content/value change for the COM surface to apply."
